$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 51-54 and 143-144: reorder/update match data (reversal of 4 rows and swap of 2 rows) ---
# Row 51
$ws.Cells.Item(51, "B").Value = 5140743
$ws.Cells.Item(51, "F").Value = 'Stal Rzeszow'
$ws.Cells.Item(51, "G").Value = 'Termalica BB Nieciecza'
$ws.Cells.Item(51, "H").Value = 2
$ws.Cells.Item(51, "I").Value = 2
$ws.Cells.Item(51, "K").Value = 3
$ws.Cells.Item(51, "L").Value = 3.3
$ws.Cells.Item(51, "M").Value = 2.2
$ws.Cells.Item(51, "N").Value = 2.9
$ws.Cells.Item(51, "O").Value = 3.3
$ws.Cells.Item(51, "P").Value = 2.25
$ws.Cells.Item(51, "Q").Value = 0.25
$ws.Cells.Item(51, "R").Value = 1.825
$ws.Cells.Item(51, "S").Value = 1.975
$ws.Cells.Item(51, "T").Value = 2.5
$ws.Cells.Item(51, "U").Value = 1.95
$ws.Cells.Item(51, "X").Value = 2.3
$ws.Cells.Item(51, "Z").Value = 0.4125
$ws.Cells.Item(51, "AA").Value = -0.5
$ws.Cells.Item(51, "AB").Value = 0.95
$ws.Cells.Item(51, "AC").Value = -1

# Row 52
$ws.Cells.Item(52, "B").Value = 5139054
$ws.Cells.Item(52, "F").Value = 'GKS Tychy 71'
$ws.Cells.Item(52, "G").Value = 'Sandecja Nowy Sacz'
$ws.Cells.Item(52, "H").Value = 2
$ws.Cells.Item(52, "I").Value = 3
$ws.Cells.Item(52, "J").Value = 'A'
$ws.Cells.Item(52, "K").Value = 2.15
$ws.Cells.Item(52, "M").Value = 3.1
$ws.Cells.Item(52, "N").Value = 2.375
$ws.Cells.Item(52, "O").Value = 3
$ws.Cells.Item(52, "P").Value = 3
$ws.Cells.Item(52, "Q").Value = -0.25
$ws.Cells.Item(52, "R").Value = 2.025
$ws.Cells.Item(52, "S").Value = 1.775
$ws.Cells.Item(52, "U").Value = 1.975
$ws.Cells.Item(52, "V").Value = 1.825
$ws.Cells.Item(52, "X").Value = -1
$ws.Cells.Item(52, "Y").Value = 2
$ws.Cells.Item(52, "Z").Value = -1
$ws.Cells.Item(52, "AA").Value = 0.7749999999999999
$ws.Cells.Item(52, "AB").Value = 0.9750000000000001
$ws.Cells.Item(52, "AC").Value = -1

# Row 53
$ws.Cells.Item(53, "B").Value = 5139053
$ws.Cells.Item(53, "F").Value = 'Chrobry Glogow'
$ws.Cells.Item(53, "G").Value = 'Zaglebie Sosnowiec'
$ws.Cells.Item(53, "H").Value = 0
$ws.Cells.Item(53, "I").Value = 0
$ws.Cells.Item(53, "J").Value = 'D'
$ws.Cells.Item(53, "K").Value = 2.45
$ws.Cells.Item(53, "M").Value = 2.55
$ws.Cells.Item(53, "N").Value = 2.7
$ws.Cells.Item(53, "O").Value = 3.2
$ws.Cells.Item(53, "P").Value = 2.375
$ws.Cells.Item(53, "Q").Value = 0
$ws.Cells.Item(53, "R").Value = 2.05
$ws.Cells.Item(53, "S").Value = 1.75
$ws.Cells.Item(53, "U").Value = 1.875
$ws.Cells.Item(53, "V").Value = 1.925
$ws.Cells.Item(53, "X").Value = 2.2
$ws.Cells.Item(53, "Y").Value = -1
$ws.Cells.Item(53, "Z").Value = 0
$ws.Cells.Item(53, "AA").Value = -0
$ws.Cells.Item(53, "AB").Value = -1
$ws.Cells.Item(53, "AC").Value = 0.925

# Row 54
$ws.Cells.Item(54, "B").Value = 5139056
$ws.Cells.Item(54, "F").Value = 'Odra Opole'
$ws.Cells.Item(54, "G").Value = 'Arka Gdynia'
$ws.Cells.Item(54, "H").Value = 1
$ws.Cells.Item(54, "I").Value = 1
$ws.Cells.Item(54, "K").Value = 3.75
$ws.Cells.Item(54, "L").Value = 3.5
$ws.Cells.Item(54, "M").Value = 1.85
$ws.Cells.Item(54, "N").Value = 3.4
$ws.Cells.Item(54, "O").Value = 3.5
$ws.Cells.Item(54, "P").Value = 1.909
$ws.Cells.Item(54, "Q").Value = 0.5
$ws.Cells.Item(54, "R").Value = 1.85
$ws.Cells.Item(54, "S").Value = 2
$ws.Cells.Item(54, "T").Value = 2.75
$ws.Cells.Item(54, "U").Value = 2
$ws.Cells.Item(54, "X").Value = 2.5
$ws.Cells.Item(54, "Z").Value = 0.8500000000000001
$ws.Cells.Item(54, "AA").Value = -1
$ws.Cells.Item(54, "AB").Value = -1
$ws.Cells.Item(54, "AC").Value = 0.8500000000000001

# Row 143
$ws.Cells.Item(143, "B").Value = 5448048
$ws.Cells.Item(143, "F").Value = 'Zaglebie Sosnowiec'
$ws.Cells.Item(143, "G").Value = 'Sandecja Nowy Sacz'
$ws.Cells.Item(143, "H").Value = 1
$ws.Cells.Item(143, "I").Value = 1
$ws.Cells.Item(143, "J").Value = 'D'
$ws.Cells.Item(143, "K").Value = 2.1
$ws.Cells.Item(143, "L").Value = 3.2
$ws.Cells.Item(143, "M").Value = 3.3
$ws.Cells.Item(143, "N").Value = 2.1
$ws.Cells.Item(143, "O").Value = 3.2
$ws.Cells.Item(143, "P").Value = 3.1
$ws.Cells.Item(143, "Q").Value = -0.25
$ws.Cells.Item(143, "R").Value = 1.875
$ws.Cells.Item(143, "S").Value = 1.925
$ws.Cells.Item(143, "T").Value = 2.25
$ws.Cells.Item(143, "U").Value = 1.85
$ws.Cells.Item(143, "V").Value = 1.95
$ws.Cells.Item(143, "X").Value = 2.2
$ws.Cells.Item(143, "Y").Value = -1
$ws.Cells.Item(143, "Z").Value = -0.5
$ws.Cells.Item(143, "AA").Value = 0.4625
$ws.Cells.Item(143, "AB").Value = -0.5
$ws.Cells.Item(143, "AC").Value = 0.475

# Row 144
$ws.Cells.Item(144, "B").Value = 5447925
$ws.Cells.Item(144, "F").Value = 'Gornik Leczna'
$ws.Cells.Item(144, "G").Value = 'Wisla Krakow'
$ws.Cells.Item(144, "H").Value = 0
$ws.Cells.Item(144, "I").Value = 3
$ws.Cells.Item(144, "J").Value = 'A'
$ws.Cells.Item(144, "K").Value = 5.5
$ws.Cells.Item(144, "L").Value = 4
$ws.Cells.Item(144, "M").Value = 1.5
$ws.Cells.Item(144, "N").Value = 4.5
$ws.Cells.Item(144, "O").Value = 4
$ws.Cells.Item(144, "P").Value = 1.615
$ws.Cells.Item(144, "Q").Value = 0.75
$ws.Cells.Item(144, "R").Value = 2.05
$ws.Cells.Item(144, "S").Value = 1.8
$ws.Cells.Item(144, "T").Value = 3
$ws.Cells.Item(144, "U").Value = 2
$ws.Cells.Item(144, "V").Value = 1.85
$ws.Cells.Item(144, "X").Value = -1
$ws.Cells.Item(144, "Y").Value = 0.615
$ws.Cells.Item(144, "Z").Value = -1
$ws.Cells.Item(144, "AA").Value = 0.8
$ws.Cells.Item(144, "AB").Value = 0
$ws.Cells.Item(144, "AC").Value = -0

# --- Rows 328-334: delete old rows 328 and 329 (data shifts up by 2) ---
$ws.Rows.Item(328).Delete()
$ws.Rows.Item(328).Delete()

# --- Apply residual corrections to shifted rows 328-332 ---
# Row 328
$ws.Cells.Item(328, "A").Value = 326
$ws.Cells.Item(328, "N").Value = 2.6
$ws.Cells.Item(328, "P").Value = 2.625
$ws.Cells.Item(328, "R").Value = 1.9
$ws.Cells.Item(328, "S").Value = 1.95

# Row 329
$ws.Cells.Item(329, "A").Value = 327

# Row 330
$ws.Cells.Item(330, "A").Value = 328

# Row 331
$ws.Cells.Item(331, "A").Value = 329
$ws.Cells.Item(331, "R").Value = 1.975
$ws.Cells.Item(331, "S").Value = 1.875
$ws.Cells.Item(331, "U").Value = 1.975
$ws.Cells.Item(331, "V").Value = 1.875

# Row 332
$ws.Cells.Item(332, "A").Value = 330

